$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Kelas" in column E, row 1, matching the bold style of the
# existing header row (A1:D1 use style index 1 / bold font)
$ws.Range("E1").Value = "Kelas"
$ws.Range("E1").Font.Bold = $true

# Update selection/active cell as in diff (activeCell D7, sqref D7)
$ws.Range("D7").Select()
